# Generate Report for Handback
# Removes the handed-back d4bd5148-... entry (row 3) from every sheet and
# refreshes the Correspond Handoff/Handback datetime stamps for the
# remaining d2f80547-... entry on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Drop all hyperlinks and remove row 3 (the d4bd5148-... entry) -------
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()

$wsZhCn.Cells.Hyperlinks.Delete()
$wsZhCn.Rows.Item(3).Delete()

$wsDeDe.Cells.Hyperlinks.Delete()
$wsDeDe.Rows.Item(3).Delete()

# --- Refresh the handoff/handback timestamps on the remaining row --------
$wsZhCn.Range("E2").Value = "2016-03-23 09:40:46"
$wsZhCn.Range("H2").Value = "2016-03-23 09:41:27"

$wsDeDe.Range("E2").Value = "2016-03-23 09:40:54"
$wsDeDe.Range("H2").Value = "2016-03-23 09:41:42"

# --- Re-create hyperlinks for the row that remains ------------------------

# Overview sheet: file-name link in column A
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/118b75cc215caeb2708e12720539304dfd57349d/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
) | Out-Null

# zh-cn sheet
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/118b75cc215caeb2708e12720539304dfd57349d/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83cec532e4322bee6af05a133d7ad35fb9b66bd0/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf"
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/d1e15ccd283f70f9c1a6774c6b9057694da20513/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75a8177be27f1a10ba82791c587e6257ce09ae44/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf"
) | Out-Null

# de-de sheet
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/118b75cc215caeb2708e12720539304dfd57349d/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c126d46bb621f71946c833821c470ed18aad4c00/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf"
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/1d2c3417ca290ffcf68d5c50baa6f71771574808/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.md"
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0baa2e8568e7dd4bc6dd74972d25fbdec151e8d4/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf",
    "",
    "",
    "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf"
) | Out-Null

$wb.Save()
